$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gets two trailing
#    spaces, then a new red run "(This is a change – Version for branch
#    alternate)" appended (split across three runs in the source diff, but
#    the visible text/formatting is what matters).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$find = $p1.Range.Find
$find.ClearFormatting()
$find.Text = "This is a Microsoft word document."
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "This is a Microsoft word document.  "
$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

$p1 = $d.Paragraphs(1)
$insertPoint = $p1.Range.End - 1
$r = $d.Range($insertPoint, $insertPoint)
$r.InsertAfter("(This is a change " + [char]0x2013 + " Version for branch alternate)")
$newEnd = $p1.Range.End - 1
$colorRange = $d.Range($insertPoint, $newEnd)
$colorRange.Font.Color = 192

Write-Output "Section 1 done"

# ---------------------------------------------------------------------------
# 2) "Crispian's Day speech ..." paragraph: tidy up the run/proofErr layout
#    - " Day speech from" gains a trailing space and the stray single-space
#      run that used to sit between it and "Shakespear's" is removed
#    - the six runs that made up " Henry V"/" "/"[Source "/"-"/" Wikipedia"/"]"
#      collapse into a single run " Henry V [Source - Wikipedia]"
# ---------------------------------------------------------------------------
$endash = [char]0x2013
$lbrk = [char]0x5B
$rbrk = [char]0x5D

$p4 = $d.Paragraphs(4)
$base = $p4.Range.Start

# remove the lone-space run between "from" and "Shakespear's"
$rSpace = $d.Range($base + 26, $base + 27)
$rSpace.Delete()

$p4 = $d.Paragraphs(4)
$find = $p4.Range.Find
$find.ClearFormatting()
$find.Text = " Day speech from"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = " Day speech from "
$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

$p4 = $d.Paragraphs(4)
$find2 = $p4.Range.Find
$find2.ClearFormatting()
$find2.Text = " Henry V " + $lbrk + "Source " + $endash + " Wikipedia" + $rbrk
$find2.Replacement.ClearFormatting()
$find2.Replacement.Text = " Henry V " + $lbrk + "Source " + $endash + " Wikipedia" + $rbrk
$find2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

Write-Output "Section 2 done"
